$wb = $excel.ActiveWorkbook

# Update selection on "Template Setup" (sheet 3) before adding the new sheet,
# so it doesn't end up marked as the tab-selected sheet.
$wsSetup = $wb.Worksheets.Item("Template Setup")
[void]$wsSetup.Range("E12").Select()

# Add the new "URL" sheet after the last existing sheet.
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$wsUrl = $wb.Worksheets.Add($null, $lastSheet)
$wsUrl.Name = "URL"

# Column A is wide enough to show the URL.
$wsUrl.Columns.Item(1).ColumnWidth = 51

# Header cell, bold + centered (reuses the existing bold font, new centered style).
$wsUrl.Range("A1").Value = "Parcel URL"
$wsUrl.Range("A1").Font.Bold = $true
$wsUrl.Range("A1").HorizontalAlignment = -4108
$wsUrl.Range("A1").VerticalAlignment = -4108

# URL cell with a real hyperlink (creates the Hyperlink style/font automatically).
$wsUrl.Range("A2").Value = "https://uat.parcelplatform.com/reporting/login.php"
$wsUrl.Hyperlinks.Add($wsUrl.Range("A2"), "https://uat.parcelplatform.com/reporting/login.php")

# Match page setup orientation of the other sheets.
$wsUrl.PageSetup.Orientation = 1

# Match the selection left on the new sheet.
[void]$wsUrl.Range("F5").Select()
